# Fruta / hortaliza, semanal
# Insert a new weekly record row at row 35 (shifts existing rows 35-54 down to 36-55)
# and populate it with this week's Chirimoya price data for Vega Monumental Concepción.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 35, pushing all following rows down by one.
$ws.Rows("35:35").Insert()

# Fill in the new row with the new weekly data.
$ws.Range("A35").Value = 11
$ws.Range("B35").Value = "Vega Monumental Concepción"
$ws.Range("C35").Value = "Bíobío"
$ws.Range("D35").Value = 45216
$ws.Range("E35").Value = 8
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100107
$ws.Range("H35").Value = "Otros"
$ws.Range("I35").Value = 100107002
$ws.Range("J35").Value = "Chirimoya"
$ws.Range("K35").Value = "Cultivar IV Región"
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 140
$ws.Range("N35").Value = 20000
$ws.Range("O35").Value = 21000
$ws.Range("P35").Value = 20429
$ws.Range("Q35").Value = "`$/bandeja 10 kilos"
$ws.Range("R35").Value = "Provincia de Limarí"
$ws.Range("S35").Value = 2043
$ws.Range("T35").Value = 10
